# Apply the "Refined metadata to be additional tab" edit:
# 1) Refresh the per-row "panel_query_time" (column F) timestamps on "data"
#    to reflect the panel-data re-scrape.
# 2) Add a new "metadata" worksheet (placed after "data") carrying summary
#    info about the scraped panel/query.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- refresh column F ("time_taken") timestamps, row by row ---
$fValues = @(
   "2021-10-05 14:22:25.648777"
  ,"2021-10-05 14:22:25.648783"
  ,"2021-10-05 14:22:25.648785"
  ,"2021-10-05 14:22:25.648787"
  ,"2021-10-05 14:22:25.648790"
  ,"2021-10-05 14:22:25.648792"
  ,"2021-10-05 14:22:25.648794"
  ,"2021-10-05 14:22:25.648795"
  ,"2021-10-05 14:22:25.648797"
  ,"2021-10-05 14:22:25.648799"
  ,"2021-10-05 14:22:25.648801"
  ,"2021-10-05 14:22:25.648803"
  ,"2021-10-05 14:22:25.648805"
  ,"2021-10-05 14:22:25.648807"
  ,"2021-10-05 14:22:25.648809"
  ,"2021-10-05 14:22:25.648811"
  ,"2021-10-05 14:22:25.648813"
  ,"2021-10-05 14:22:25.648815"
  ,"2021-10-05 14:22:25.648817"
  ,"2021-10-05 14:22:25.648819"
  ,"2021-10-05 14:22:25.648821"
  ,"2021-10-05 14:22:25.648823"
  ,"2021-10-05 14:22:25.648825"
  ,"2021-10-05 14:22:25.648827"
  ,"2021-10-05 14:22:25.648829"
  ,"2021-10-05 14:22:25.648831"
  ,"2021-10-05 14:22:25.648833"
  ,"2021-10-05 14:22:25.648835"
  ,"2021-10-05 14:22:25.648837"
  ,"2021-10-05 14:22:25.648839"
  ,"2021-10-05 14:22:25.648841"
  ,"2021-10-05 14:22:25.648843"
  ,"2021-10-05 14:22:25.648845"
  ,"2021-10-05 14:22:25.648847"
  ,"2021-10-05 14:22:25.648849"
  ,"2021-10-05 14:22:25.648850"
  ,"2021-10-05 14:22:25.648852"
  ,"2021-10-05 14:22:25.648854"
  ,"2021-10-05 14:22:25.648856"
  ,"2021-10-05 14:22:25.648858"
  ,"2021-10-05 14:22:25.648860"
  ,"2021-10-05 14:22:25.648862"
  ,"2021-10-05 14:22:25.648864"
  ,"2021-10-05 14:22:25.648866"
  ,"2021-10-05 14:22:25.648868"
  ,"2021-10-05 14:22:25.648870"
  ,"2021-10-05 14:22:25.648872"
  ,"2021-10-05 14:22:25.648874"
  ,"2021-10-05 14:22:25.648876"
  ,"2021-10-05 14:22:25.648878"
  ,"2021-10-05 14:22:25.648880"
  ,"2021-10-05 14:22:25.648882"
  ,"2021-10-05 14:22:25.648884"
  ,"2021-10-05 14:22:25.648886"
  ,"2021-10-05 14:22:25.648888"
  ,"2021-10-05 14:22:25.648890"
  ,"2021-10-05 14:22:25.648892"
  ,"2021-10-05 14:22:25.648894"
  ,"2021-10-05 14:22:25.648895"
  ,"2021-10-05 14:22:25.648897"
  ,"2021-10-05 14:22:25.648899"
  ,"2021-10-05 14:22:25.648901"
  ,"2021-10-05 14:22:25.648903"
  ,"2021-10-05 14:22:25.648905"
  ,"2021-10-05 14:22:25.648908"
  ,"2021-10-05 14:22:25.648910"
  ,"2021-10-05 14:22:25.648912"
  ,"2021-10-05 14:22:25.648914"
  ,"2021-10-05 14:22:25.648916"
  ,"2021-10-05 14:22:25.648918"
  ,"2021-10-05 14:22:25.648920"
  ,"2021-10-05 14:22:25.648922"
  ,"2021-10-05 14:22:25.648923"
  ,"2021-10-05 14:22:25.648925"
  ,"2021-10-05 14:22:25.648927"
  ,"2021-10-05 14:22:25.648929"
  ,"2021-10-05 14:22:25.648932"
  ,"2021-10-05 14:22:25.648935"
  ,"2021-10-05 14:22:25.648937"
  ,"2021-10-05 14:22:25.648939"
  ,"2021-10-05 14:22:25.648940"
  ,"2021-10-05 14:22:25.648942"
  ,"2021-10-05 14:22:25.648945"
  ,"2021-10-05 14:22:25.648947"
  ,"2021-10-05 14:22:25.648949"
  ,"2021-10-05 14:22:25.648950"
  ,"2021-10-05 14:22:25.648952"
  ,"2021-10-05 14:22:25.648954"
  ,"2021-10-05 14:22:25.648956"
  ,"2021-10-05 14:22:25.648958"
  ,"2021-10-05 14:22:25.648960"
  ,"2021-10-05 14:22:25.648962"
  ,"2021-10-05 14:22:25.648965"
  ,"2021-10-05 14:22:25.648967"
  ,"2021-10-05 14:22:25.648969"
  ,"2021-10-05 14:22:25.648971"
  ,"2021-10-05 14:22:25.648973"
  ,"2021-10-05 14:22:25.648975"
  ,"2021-10-05 14:22:25.648977"
  ,"2021-10-05 14:22:25.648979"
  ,"2021-10-05 14:22:25.648980"
  ,"2021-10-05 14:22:25.648982"
  ,"2021-10-05 14:22:25.648984"
  ,"2021-10-05 14:22:25.648986"
  ,"2021-10-05 14:22:25.648988"
  ,"2021-10-05 14:22:25.648990"
  ,"2021-10-05 14:22:25.648992"
  ,"2021-10-05 14:22:25.648993"
  ,"2021-10-05 14:22:25.648997"
  ,"2021-10-05 14:22:25.648999"
  ,"2021-10-05 14:22:25.649001"
  ,"2021-10-05 14:22:25.649003"
  ,"2021-10-05 14:22:25.649005"
  ,"2021-10-05 14:22:25.649007"
  ,"2021-10-05 14:22:25.649009"
  ,"2021-10-05 14:22:25.649011"
  ,"2021-10-05 14:22:25.649013"
  ,"2021-10-05 14:22:25.649015"
  ,"2021-10-05 14:22:25.649016"
  ,"2021-10-05 14:22:25.649018"
  ,"2021-10-05 14:22:25.649020"
  ,"2021-10-05 14:22:25.649022"
  ,"2021-10-05 14:22:25.649024"
  ,"2021-10-05 14:22:25.649025"
  ,"2021-10-05 14:22:25.649027"
  ,"2021-10-05 14:22:25.649029"
  ,"2021-10-05 14:22:25.649031"
  ,"2021-10-05 14:22:25.649033"
  ,"2021-10-05 14:22:25.649036"
  ,"2021-10-05 14:22:25.649038"
  ,"2021-10-05 14:22:25.649040"
  ,"2021-10-05 14:22:25.649042"
  ,"2021-10-05 14:22:25.649044"
  ,"2021-10-05 14:22:25.649046"
  ,"2021-10-05 14:22:25.649048"
  ,"2021-10-05 14:22:25.649050"
  ,"2021-10-05 14:22:25.649051"
  ,"2021-10-05 14:22:25.649053"
  ,"2021-10-05 14:22:25.649055"
  ,"2021-10-05 14:22:25.649057"
  ,"2021-10-05 14:22:25.649059"
  ,"2021-10-05 14:22:25.649061"
  ,"2021-10-05 14:22:25.649063"
  ,"2021-10-05 14:22:25.649065"
  ,"2021-10-05 14:22:25.649067"
  ,"2021-10-05 14:22:25.649068"
  ,"2021-10-05 14:22:25.649070"
  ,"2021-10-05 14:22:25.649072"
  ,"2021-10-05 14:22:25.649074"
  ,"2021-10-05 14:22:25.649076"
  ,"2021-10-05 14:22:25.649078"
  ,"2021-10-05 14:22:25.649080"
  ,"2021-10-05 14:22:25.649082"
  ,"2021-10-05 14:22:25.649084"
  ,"2021-10-05 14:22:25.649086"
  ,"2021-10-05 14:22:25.649088"
  ,"2021-10-05 14:22:25.649090"
  ,"2021-10-05 14:22:25.649091"
  ,"2021-10-05 14:22:25.649093"
  ,"2021-10-05 14:22:25.649095"
  ,"2021-10-05 14:22:25.649097"
  ,"2021-10-05 14:22:25.649099"
  ,"2021-10-05 14:22:25.649101"
  ,"2021-10-05 14:22:25.649103"
  ,"2021-10-05 14:22:25.649105"
  ,"2021-10-05 14:22:25.649107"
  ,"2021-10-05 14:22:25.649109"
  ,"2021-10-05 14:22:25.649111"
  ,"2021-10-05 14:22:25.649112"
  ,"2021-10-05 14:22:25.649114"
  ,"2021-10-05 14:22:25.649116"
  ,"2021-10-05 14:22:25.649118"
  ,"2021-10-05 14:22:25.649121"
  ,"2021-10-05 14:22:25.649123"
  ,"2021-10-05 14:22:25.649125"
  ,"2021-10-05 14:22:25.649127"
  ,"2021-10-05 14:22:25.649129"
  ,"2021-10-05 14:22:25.649131"
  ,"2021-10-05 14:22:25.649133"
  ,"2021-10-05 14:22:25.649135"
  ,"2021-10-05 14:22:25.649136"
  ,"2021-10-05 14:22:25.649138"
  ,"2021-10-05 14:22:25.649140"
  ,"2021-10-05 14:22:25.649142"
  ,"2021-10-05 14:22:25.649144"
  ,"2021-10-05 14:22:25.649145"
  ,"2021-10-05 14:22:25.649147"
  ,"2021-10-05 14:22:25.649149"
  ,"2021-10-05 14:22:25.649151"
  ,"2021-10-05 14:22:25.649153"
  ,"2021-10-05 14:22:25.649155"
  ,"2021-10-05 14:22:25.649157"
  ,"2021-10-05 14:22:25.649159"
  ,"2021-10-05 14:22:25.649160"
  ,"2021-10-05 14:22:25.649162"
  ,"2021-10-05 14:22:25.649164"
  ,"2021-10-05 14:22:25.649166"
  ,"2021-10-05 14:22:25.649168"
  ,"2021-10-05 14:22:25.649170"
  ,"2021-10-05 14:22:25.649172"
  ,"2021-10-05 14:22:25.649174"
  ,"2021-10-05 14:22:25.649176"
  ,"2021-10-05 14:22:25.649178"
  ,"2021-10-05 14:22:25.649180"
)
for ($i = 0; $i -lt $fValues.Length; $i++) {
    $data.Cells.Item($i + 2, 6).Value = $fValues[$i]
}

# --- add the "metadata" sheet, placed right after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$meta = $wb.Worksheets.Add($null, $lastSheet)
$meta.Name = "metadata"

# Reuse "data"'s header style (bold font + thin border + center/top
# alignment) for the metadata header row and its index cell, by copying
# formats from the existing sheet instead of re-deriving the same style.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Rare multisystem ciliopathy disorders"
$meta.Range("C2").Value = 150
# Force "1.146" to be stored as text (not coerced to the number 1.146) with
# a leading quote-prefix, then immediately strip the quote-prefix styling
# back off by re-pasting a plain/default format over the cell - the stored
# value stays text, but the cell's format reverts to the sheet's default.
$meta.Range("D2").Value = "'1.146"
$data.Range("H1").Copy()
$meta.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$meta.Range("E2").Value = "2021-09-06T13:40:13.491630Z"
$meta.Range("F2").Value = "2021-10-05 14:22:25.646341"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/150/?format=json"

$data.Activate()
Write-Output "metadata sheet added; timestamps refreshed"
